$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed A2 with the date number format (mm-dd-yy, numFmtId 14) and a value,
# then copy that format to every other date cell so they all share one
# style index instead of each creating its own duplicate style.
$ws.Range("A2").NumberFormat = "mm-dd-yy"
$ws.Range("A2").Value = (Get-Date -Year 2022 -Month 4 -Day 23 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Range("A2").Copy()
$ws.Range("A3:A8").PasteSpecial(-4122)
$ws.Range("D2:D8").PasteSpecial(-4122)

# Row 2 - Stack push / null
$ws.Range("B2").Value = "Eric"
$ws.Range("C2").Value = "Stack function push accepts null as valid inputs"
$ws.Range("D2").Value = (Get-Date -Year 2022 -Month 4 -Day 20 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Range("E2").Value = "Duy"
$ws.Range("F2").Value = "Added an if statement to catch null and restart the function"

# Row 3 - Stack push / empty string
$ws.Range("A3").Value = (Get-Date -Year 2022 -Month 4 -Day 20 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Range("B3").Value = "Duy"
$ws.Range("C3").Value = "Stack function push accepts empty string as valid inputs"
$ws.Range("D3").Value = (Get-Date -Year 2022 -Month 4 -Day 20 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Range("E3").Value = "Duy"
$ws.Range("F3").Value = "Added an if statement to catch empty string and restart the function"

# Row 4 - Queue enqueue / empty string
$ws.Range("A4").Value = (Get-Date -Year 2022 -Month 4 -Day 20 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Range("B4").Value = "Duy"
$ws.Range("C4").Value = "Queue function enqueue accepts empty strings as valid input"
$ws.Range("D4").Value = (Get-Date -Year 2022 -Month 4 -Day 21 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Range("E4").Value = "Aidan"
$ws.Range("F4").Value = "Added if statement to return to main function page after user enters empty string"

# Row 5 - Queue enqueue / cancel
$ws.Range("A5").Value = (Get-Date -Year 2022 -Month 4 -Day 21 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Range("B5").Value = "Aidan "
$ws.Range("C5").Value = "Queue function enqueue accepts cancel input from user"
$ws.Range("D5").Value = (Get-Date -Year 2022 -Month 4 -Day 21 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Range("E5").Value = "Aidan"
$ws.Range("F5").Value = "Added if statement to catch the cancel input which gives null and returns the user into the main page of queue"

# Row 6 - Linked List insert / empty string
$ws.Range("A6").Value = (Get-Date -Year 2022 -Month 4 -Day 21 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Range("B6").Value = "Aidan "
$ws.Range("C6").Value = "Linked List insert accepts empty string as a valid input"
$ws.Range("D6").Value = (Get-Date -Year 2022 -Month 4 -Day 22 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Range("E6").Value = "Eric"
$ws.Range("F6").Value = "Added another if statement to catch the empty string to return the user back to main page of linked list so they can input again"

# Row 7 - Linked List insert / cancel
$ws.Range("A7").Value = (Get-Date -Year 2022 -Month 4 -Day 22 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Range("B7").Value = "Eric"
$ws.Range("C7").Value = "Linked List insert accepts cancel from user as a valid input"
$ws.Range("D7").Value = (Get-Date -Year 2022 -Month 4 -Day 22 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Range("E7").Value = "Eric "
$ws.Range("F7").Value = "Added if statement to make sure that cancel is not accepted as a input and returns the user back to the main page of linked list to input a valid input"

# Row 8 - BST nodes (no A8/B8 values; clear the date format seeded there)
$ws.Range("A8").Clear()
$ws.Range("C8").Value = "BST nodes does not appear as expected"
$ws.Range("D8").Value = (Get-Date -Year 2022 -Month 4 -Day 23 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Range("E8").Value = "Landen"
$ws.Range("F8").Value = "Switch the browser to run the program"

$ws.Range("C8:F8").Select()
